# Re-save the "Authors" (column E) values for the data rows.
#
# The source data pipeline re-emitted the Authors column text for every
# reference row (rows 2-33, skipping the three rows whose Authors value is
# the empty-list placeholder "[]"). The re-emitted text is identical in
# content (same names/emails/numbers) but every comma separator gained one
# extra padding space, e.g. ",                " -> ",                 ".
# This duplicates/repoints the shared-string table the same way the
# original commit did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    $current = $cell.Value2
    if ($current -eq "[]") {
        continue
    }
    $updated = $current -replace ',( +)', ', $1'
    $cell.Value2 = $updated
}
